# Added windows settings, panel settings screens, added animation between screens.
#
# Updates the "Translation" sheet of the TouchGFX texts workbook:
#   - E7  : "0"       -> "00"
#   - E12 : "Monday"  -> "Wednesday"
#   - Rows 13-24: fill in the previously-empty placeholder rows with new
#     Text ID / Typography Name / Alignment / GB / Direction entries for
#     the new Window Settings / Zones / Scenes / Panel Settings screens.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

# --- E7: "0" -> "00" -----------------------------------------------------
# A bare "00" would be auto-converted to the number 0 by the smart-entry
# parser, so type it in as quote-prefixed text (keeps it a string) and then
# reset the cell style back to Normal so no stray per-cell style sticks
# around once the value has been committed.
$ws.Range("E7").Value = "'00"
$ws.Range("E7").Style = "Normal"

# --- E12: "Monday" -> "Wednesday" ----------------------------------------
$ws.Range("E12").Value = "Wednesday"

# --- New rows 13-24 --------------------------------------------------------
$newRows = @(
    @{ Row = 13; B = "SingleUseId10"; C = "Default"; D = "Left";   E = "Window Settings";    F = "LTR" },
    @{ Row = 14; B = "SingleUseId11"; C = "Default"; D = "Center"; E = "Zones";              F = "LTR" },
    @{ Row = 15; B = "SingleUseId12"; C = "Default"; D = "Center"; E = "Scenes";             F = "LTR" },
    @{ Row = 16; B = "SingleUseId26"; C = "Default"; D = "Left";   E = "Checkerboard";       F = "LTR" },
    @{ Row = 17; B = "SingleUseId21"; C = "Default"; D = "Left";   E = "Geometric Pattern";  F = "LTR" },
    @{ Row = 18; B = "SingleUseId23"; C = "Default"; D = "Left";   E = '"Go Navy" Text';     F = "LTR" },
    @{ Row = 19; B = "SingleUseId25"; C = "Default"; D = "Left";   E = "Sine Wave";          F = "LTR" },
    @{ Row = 20; B = "SingleUseId19"; C = "Default"; D = "Left";   E = "Christmas Tree";     F = "LTR" },
    @{ Row = 21; B = "SingleUseId27"; C = "Default"; D = "Left";   E = "Auto Tint";          F = "LTR" },
    @{ Row = 22; B = "SingleUseId28"; C = "Default"; D = "Left";   E = "Manual Tint";        F = "LTR" },
    @{ Row = 23; B = "SingleUseId29"; C = "Default"; D = "Left";   E = "Preset Scenes";      F = "LTR" },
    @{ Row = 24; B = "SingleUseId30"; C = "Default"; D = "Center"; E = "Panel Settings";     F = "LTR" }
)

foreach ($entry in $newRows) {
    $r = $entry.Row
    $ws.Range("B$r").Value = $entry.B
    $ws.Range("C$r").Value = $entry.C
    $ws.Range("D$r").Value = $entry.D
    $ws.Range("E$r").Value = $entry.E
    $ws.Range("F$r").Value = $entry.F
    # Writing into a previously-blank cell stamps the inherited column style
    # onto the cell explicitly; reset to Normal so the cells stay styleless,
    # matching the rest of the table's data rows.
    $ws.Range("B$r`:F$r").Style = "Normal"
}
